$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 2.9
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3.5
$ws.Range("L3").Value = 4
$ws.Range("Y3").Value = 2.38
$ws.Range("Z3").Value = 1.53
$ws.Range("AA3").Value = 5.5
$ws.Range("AB3").Value = 10
$ws.Range("AC3").Value = 11
$ws.Range("AD3").Value = 26
$ws.Range("AF3").Value = 41
$ws.Range("AJ3").Value = 101
$ws.Range("AL3").Value = 6.5
$ws.Range("AM3").Value = 13
$ws.Range("AN3").Value = 13
$ws.Range("AO3").Value = 34
$ws.Range("AP3").Value = 34
$ws.Range("Q4").Value = 2.4
$ws.Range("R4").Value = 1.53
$ws.Range("S4").Value = 4.1
$ws.Range("T4").Value = 1.24
$ws.Range("U4").Value = 4.5
$ws.Range("V4").Value = 1.18
$ws.Range("AR4").Value = 1.85
$ws.Range("AS4").Value = 2
$ws.Range("S6").Value = 4.6
$ws.Range("T6").Value = 1.2
$ws.Range("AR6").Value = 2.1
$ws.Range("AS6").Value = 1.78
$ws.Range("L7").Value = 2.1
$ws.Range("AF7").Value = 67
$ws.Range("AM7").Value = 6
$ws.Range("Z8").Value = 1.47
$ws.Range("G9").Value = 1.44
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 8.5
$ws.Range("J9").Value = 2.05
$ws.Range("K9").Value = 2.1
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 3.35
$ws.Range("T9").Value = 1.32
$ws.Range("W9").Value = 1.5
$ws.Range("X9").Value = 2.5
$ws.Range("Z9").Value = 1.47
$ws.Range("AD9").Value = 9
$ws.Range("AL9").Value = 15
$ws.Range("AN9").Value = 26
$ws.Range("AP9").Value = 67
$ws.Range("AR9").Value = 1.67
$ws.Range("AS9").Value = 2.17
$ws.Range("Q10").Value = 2.6
$ws.Range("R10").Value = 1.48
$ws.Range("S10").Value = 4.1
$ws.Range("T10").Value = 1.23
$ws.Range("Z10").Value = 1.63
$ws.Range("AR10").Value = 1.95
$ws.Range("AS10").Value = 1.9
$ws.Range("G12").Value = 3.15
$ws.Range("H12").Value = 2.72
$ws.Range("I12").Value = 2.55
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 1.88
$ws.Range("L12").Value = 3.15
$ws.Range("M12").Value = 1.12
$ws.Range("N12").Value = 5.3
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.42
$ws.Range("Q12").Value = 2.42
$ws.Range("R12").Value = 1.5
$ws.Range("U12").Value = 4.35
$ws.Range("V12").Value = 1.18
$ws.Range("W12").Value = 1.53
$ws.Range("X12").Value = 2.35
$ws.Range("Y12").Value = 1.98
$ws.Range("Z12").Value = 1.75
$ws.Range("AA12").Value = 7.3
$ws.Range("AC12").Value = 11.25
$ws.Range("AD12").Value = 45
$ws.Range("AE12").Value = 35
$ws.Range("AF12").Value = 50
$ws.Range("AG12").Value = 5.3
$ws.Range("AH12").Value = 5.4
$ws.Range("AI12").Value = 15.5
$ws.Range("AJ12").Value = 90
$ws.Range("AK12").Value = 900
$ws.Range("AL12").Value = 6.5
$ws.Range("AM12").Value = 11.75
$ws.Range("AN12").Value = 9.75
$ws.Range("AO12").Value = 29
$ws.Range("AP12").Value = 25
$ws.Range("AQ12").Value = 40
$ws.Range("G13").Value = 2.15
$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = 2.82
$ws.Range("K13").Value = 1.87
$ws.Range("L13").Value = 4.6
$ws.Range("M13").Value = 1.14
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 1.55
$ws.Range("P13").Value = 2.3
$ws.Range("Q13").Value = 2.62
$ws.Range("R13").Value = 1.44
$ws.Range("U13").Value = 4.7
$ws.Range("V13").Value = 1.15
$ws.Range("W13").Value = 1.57
$ws.Range("X13").Value = 2.27
$ws.Range("Y13").Value = 2.12
$ws.Range("Z13").Value = 1.65
$ws.Range("AA13").Value = 5.4
$ws.Range("AC13").Value = 9.25
$ws.Range("AD13").Value = 21
$ws.Range("AE13").Value = 22
$ws.Range("AF13").Value = 40
$ws.Range("AG13").Value = 5
$ws.Range("AI13").Value = 18
$ws.Range("AJ13").Value = 120
$ws.Range("AL13").Value = 8.25
$ws.Range("AN13").Value = 14
$ws.Range("AO13").Value = 70
$ws.Range("AP13").Value = 50
$ws.Range("AQ13").Value = 65
$ws.Range("G14").Value = 1.32
$ws.Range("I14").Value = 10
$ws.Range("AL14").Value = 19.5
$ws.Range("G16").Value = 1.62
$ws.Range("H16").Value = 3.7
$ws.Range("I16").Value = 5.75
$ws.Range("W16").Value = 1.44
$ws.Range("X16").Value = 2.63
$ws.Range("AI16").Value = 21
$ws.Range("AM16").Value = 29
$ws.Range("AN16").Value = 19
$ws.Range("AO16").Value = 67
$ws.Range("AP16").Value = 51
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.5
$ws.Range("AR17").Value = 1.9
$ws.Range("AS17").Value = 1.95
$ws.Range("G19").Value = 2.4
$ws.Range("L19").Value = 3.5
$ws.Range("AG19").Value = 9.5
$ws.Range("M20").Value = 1.06
$ws.Range("N20").Value = 10
$ws.Range("O20").Value = 1.3
$ws.Range("P20").Value = 3.5
$ws.Range("Q20").Value = 1.98
$ws.Range("R20").Value = 1.88
$ws.Range("O21").Value = 1.33
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 2.05
$ws.Range("R21").Value = 1.8
$ws.Range("G22").Value = 2.3
$ws.Range("I22").Value = 2.8
$ws.Range("L22").Value = 3.4
$ws.Range("AA22").Value = 9.5
$ws.Range("AB22").Value = 12
$ws.Range("AE22").Value = 17
$ws.Range("AN22").Value = 11
$ws.Range("G23").Value = 1.9
$ws.Range("I23").Value = 4.33
$ws.Range("J23").Value = 2.63
$ws.Range("O23").Value = 1.5
$ws.Range("P23").Value = 2.5
$ws.Range("Q23").Value = 2.5
$ws.Range("R23").Value = 1.5
$ws.Range("U23").Value = 5
$ws.Range("V23").Value = 1.17
$ws.Range("AD23").Value = 15
$ws.Range("AL23").Value = 9.5
$ws.Range("AM23").Value = 21
$ws.Range("AO23").Value = 51
$ws.Range("AR23").Value = 1.88
$ws.Range("AS23").Value = 1.93
$ws.Range("O24").Value = 1.3
$ws.Range("P24").Value = 3.4
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = 1.8
$ws.Range("U24").Value = 3.5
$ws.Range("V24").Value = 1.29
$ws.Range("M25").Value = 1.07
$ws.Range("N25").Value = 9
$ws.Range("O25").Value = 1.33
$ws.Range("P25").Value = 3.25
$ws.Range("Q25").Value = 2.1
$ws.Range("R25").Value = 1.7
$ws.Range("U25").Value = 3.75
$ws.Range("V25").Value = 1.25
$ws.Range("Y25").Value = 1.8
$ws.Range("Z25").Value = 1.91
$ws.Range("AF25").Value = 34
$ws.Range("AJ25").Value = 51
$ws.Range("AK25").Value = 251
$ws.Range("AL25").Value = 8.5
$ws.Range("AP25").Value = 23
